# Regenerate save_data column G ("K") values for rows 2-25.
# This mirrors the upstream commit that recalculated K (strikeouts)
# instead of the previous Strike# figure, and rewrote the resulting
# s_vals into the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 0
    4  = 1
    5  = 2
    6  = 2
    7  = 2
    8  = 2
    9  = 1
    10 = 3
    11 = 1
    12 = 1
    13 = 1
    14 = 0
    15 = 1
    16 = 1
    17 = 1
    18 = 1
    19 = 2
    20 = 2
    21 = 2
    22 = 2
    23 = 1
    24 = 1
    25 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
